# M14 Froze Encoder 1234
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update numeric score column (C) values
$ws.Range("C2").Value  = 18
$ws.Range("C3").Value  = 10
$ws.Range("C5").Value  = 21
$ws.Range("C6").Value  = 13
$ws.Range("C7").Value  = 22
$ws.Range("C9").Value  = 11
$ws.Range("C11").Value = 13
$ws.Range("C12").Value = 14
$ws.Range("C14").Value = 14
$ws.Range("C15").Value = 19
$ws.Range("C16").Value = 15
$ws.Range("C17").Value = 20
$ws.Range("C18").Value = 18

# Update hypothesis text column (B) values
$ws.Range("B10").Value = "<hin>"
$ws.Range("B18").Value = "<unifonm>"
